$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the existing header style/format (border, bold, centered)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill in the new "Save" column values for the data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
